# Update ODI match counts for a set of active players as per additional scraping.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 26
$ws.Range("C5").Value = 140
$ws.Range("C14").Value = 9
$ws.Range("C15").Value = 36
$ws.Range("C17").Value = 7
$ws.Range("C20").Value = 50
$ws.Range("C21").Value = 155
$ws.Range("C23").Value = 45
$ws.Range("C24").Value = 21
$ws.Range("C29").Value = 89
$ws.Range("C32").Value = 44
$ws.Range("C36").Value = 45
